# Update machine_spec.xlsx with MEC prod data.
# Re-orders the header/value columns (descr/lang_code swap position) and
# appends the audit columns (cr_by, cr_dtimes, upd_by, upd_dtimes,
# is_deleted, del_dtimes) plus a brand-new data row for the Android
# machine-spec record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 1 - headers (existing A1:H1 reshuffled, I1 "is_active" unchanged,
# J1:O1 brand new audit-column headers)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "brand"
$ws.Range("D1").Value = "model"
$ws.Range("E1").Value = "mtyp_code"
$ws.Range("F1").Value = "min_driver_ver"
$ws.Range("G1").Value = "descr"
$ws.Range("H1").Value = "lang_code"
$ws.Range("J1").Value = "cr_by"
$ws.Range("K1").Value = "cr_dtimes"
$ws.Range("L1").Value = "upd_by"
$ws.Range("M1").Value = "upd_dtimes"
$ws.Range("N1").Value = "is_deleted"
$ws.Range("O1").Value = "del_dtimes"

# ---------------------------------------------------------------------
# Row 2 - existing "Resident Virtual Machine" record, columns reshuffled
# the same way as the header row, plus the new audit values
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "RESIDENT-1"
$ws.Range("B2").Value = "Resident Virtual Machine"
$ws.Range("C2").Value = "Unkown"
$ws.Range("D2").Value = "Unknown"
$ws.Range("E2").Value = "RESIDENT-REG"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "Resident Virtual Machine"
$ws.Range("H2").Value = "eng"
$ws.Range("I2").Value = $true
$ws.Range("J2").Value = "rediet"
$ws.Range("K2").Value = 44776.354369108798
$ws.Range("K2").NumberFormat = "mm:ss.0"
$ws.Range("N2").Value = $false

# ---------------------------------------------------------------------
# Row 3 - brand new "Android" machine-spec record
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "3ce7764d-42c3-4851-9d88-eb40e755b09b"
$ws.Range("B3").Value = "Android"
$ws.Range("C3").Value = "Android"
$ws.Range("D3").Value = "Android"
$ws.Range("E3").Value = "ANDROID"
$ws.Range("F3").Value = "Android"
$ws.Range("G3").Value = "Android"
$ws.Range("H3").Value = "eng"
$ws.Range("I3").Value = $true
$ws.Range("J3").Value = 1103
$ws.Range("K3").Value = 44825.475885162035
$ws.Range("K3").NumberFormat = "mm:ss.0"
$ws.Range("L3").Value = 1103
$ws.Range("M3").Value = 44825.476009502316
$ws.Range("M3").NumberFormat = "mm:ss.0"
$ws.Range("N3").Value = $false

# ---------------------------------------------------------------------
# Misc view state tweak that came along with the data refresh
# ---------------------------------------------------------------------
$ws.Range("C8").Select() | Out-Null
